$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Val)
    if ($Val -match '^\s*[-+]?\d+(\.\d+)?\s*$' -and $Val -notmatch '\s') {
        $Range.Value = "'" + $Val
    } else {
        $Range.Value = $Val
    }
}

Set-TextValue $ws.Range("D2") '26.874.34'
Set-TextValue $ws.Range("E2") '  +0.70%  '
Set-TextValue $ws.Range("D3") '1.642.34'
Set-TextValue $ws.Range("E3") '  +0.63%  '
Set-TextValue $ws.Range("E4") '  -0.52%  '
Set-TextValue $ws.Range("D5") '216.56'
Set-TextValue $ws.Range("E5") '  -0.63%  '
Set-TextValue $ws.Range("E6") '  +1.96%  '
Set-TextValue $ws.Range("E7") '  -0.55%  '
Set-TextValue $ws.Range("E8") '  +1.94%  '
Set-TextValue $ws.Range("D9") '0.0621'
Set-TextValue $ws.Range("E9") '  +0.19%  '
Set-TextValue $ws.Range("D10") '19.82'
Set-TextValue $ws.Range("E10") '  +4.41%  '
Set-TextValue $ws.Range("D12") '1.872.13'
Set-TextValue $ws.Range("E12") '  +0.64%  '
Set-TextValue $ws.Range("D13") '1.654.05'
Set-TextValue $ws.Range("E13") '  +1.37%  '
Set-TextValue $ws.Range("E14") '  +0.57%  '
Set-TextValue $ws.Range("E15") '  +1.47%  '
Set-TextValue $ws.Range("D16") '66.38'
Set-TextValue $ws.Range("E16") '  +3.79%  '
Set-TextValue $ws.Range("D17") '26.881.06'
Set-TextValue $ws.Range("E17") '  +0.79%  '
Set-TextValue $ws.Range("D18") '0.0₃0728'
Set-TextValue $ws.Range("E18") '  +0.90%  '
Set-TextValue $ws.Range("D19") '219.44'
Set-TextValue $ws.Range("E19") '  +3.95%  '
Set-TextValue $ws.Range("E20") '  -0.65%  '
Set-TextValue $ws.Range("B21") 'Uniswap'
Set-TextValue $ws.Range("C21") 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range("D21") '4.38'
Set-TextValue $ws.Range("E21") '  +1.89%  '
Set-TextValue $ws.Range("B22") 'Chainlink'
Set-TextValue $ws.Range("C22") 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D22") '6.63'
Set-TextValue $ws.Range("E22") '  +7.50%  '
Set-TextValue $ws.Range("D23") '2.42'
Set-TextValue $ws.Range("E23") '  +3.70%  '
Set-TextValue $ws.Range("D24") '9.17'
Set-TextValue $ws.Range("E24") '  +0.26%  '
Set-TextValue $ws.Range("D25") '145.83'
Set-TextValue $ws.Range("E25") '  -0.46%  '
Set-TextValue $ws.Range("E26") '  -0.71%  '
Set-TextValue $ws.Range("E27") '  +5.64%  '
Set-TextValue $ws.Range("E28") '  +1.63%  '
Set-TextValue $ws.Range("D29") '15.81'
Set-TextValue $ws.Range("E29") '  +2.05%  '
Set-TextValue $ws.Range("D30") '0.0508'
Set-TextValue $ws.Range("E30") '  +1.40%  '
Set-TextValue $ws.Range("E31") '  -0.56%  '
Set-TextValue $ws.Range("E32") '  -0.54%  '
Set-TextValue $ws.Range("E33") '  +2.05%  '
Set-TextValue $ws.Range("E34") '  +3.03%  '
Set-TextValue $ws.Range("E35") '  +0.04%  '
Set-TextValue $ws.Range("D36") '1.244.39'
Set-TextValue $ws.Range("E36") '  -1.20%  '
Set-TextValue $ws.Range("E37") '  +0.92%  '
Set-TextValue $ws.Range("E38") '  +3.24%  '
Set-TextValue $ws.Range("D39") '0.832'
Set-TextValue $ws.Range("E39") '  +3.90%  '
Set-TextValue $ws.Range("E40") '  -0.61%  '
Set-TextValue $ws.Range("D41") '0.807'
Set-TextValue $ws.Range("E41") '  +1.29%  '
Set-TextValue $ws.Range("E42") '  +2.50%  '
Set-TextValue $ws.Range("D43") '1.783.79'
Set-TextValue $ws.Range("E43") '  +0.74%  '
Set-TextValue $ws.Range("E44") '  -3.05%  '
Set-TextValue $ws.Range("D45") '60.74'
Set-TextValue $ws.Range("E45") '  +1.58%  '
Set-TextValue $ws.Range("D46") '91.48'
Set-TextValue $ws.Range("E46") '  +0.62%  '
Set-TextValue $ws.Range("E47") '  +0.80%  '
Set-TextValue $ws.Range("B48") 'Cronos'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D48") '0.0514'
Set-TextValue $ws.Range("E48") '  -0.21%  '
Set-TextValue $ws.Range("B49") 'Algorand'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D49") '0.0974'
Set-TextValue $ws.Range("E49") '  +2.22%  '
Set-TextValue $ws.Range("B50") 'EnergySwap'
Set-TextValue $ws.Range("C50") 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D50") '7.59'
Set-TextValue $ws.Range("E50") '  +1.81%  '
Set-TextValue $ws.Range("B51") 'Mantle'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range("D51") '0.405'
Set-TextValue $ws.Range("E51") '  -0.24%  '
